# Weekly price-sheet update: a new week's record (for "Locoto", quality
# "Primera") is inserted as the new row 29, pushing every existing row
# from 29..118 down by one (to 30..119). All the shifted rows keep their
# original values untouched; only the brand-new row 29 gets fresh data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 29 - this shifts rows 29:118 down to
# 30:119 and carries formatting (e.g. the date style on column D) down
# with them, same as Excel's native "Insert Sheet Rows" command.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with this week's record.
$ws.Range("A29").Value = 1
$ws.Range("B29").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C29").Value = "Arica y Parinacota"
$ws.Range("D29").Value = 44764
$ws.Range("E29").Value = 15
$ws.Range("F29").Value = 100112042
$ws.Range("G29").Value = "Locoto"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 24000
$ws.Range("L29").Value = 25000
$ws.Range("M29").Value = 24500
$ws.Range("N29").Value = "`$/caja 20 kilos"
$ws.Range("O29").Value = "Región de Arica y Parinacota"
$ws.Range("P29").Value = 1225
$ws.Range("Q29").Value = 20
$ws.Range("R29").Value = "Hortaliza"
